# Predicted_LaLiga2025_26_table_matchday_10.xlsx
# Add extra prediction columns (WIN / TOP4 / TOP5 / TOP6 / RELEGATION) between
# the existing "Team" and "ExpPoints" columns, shifting ExpPoints from C to H,
# refresh the ExpPoints values and re-order a few teams.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ------------------------------------------------------------
# Copy the header style that currently lives on C1 onto the new header cells
# before moving "ExpPoints" out to column H.
$ws.Range("C1").Copy($ws.Range("D1:H1"))

$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

# --- Data rows ---------------------------------------------------------------
# Clear the old ExpPoints numbers out of column C (now "WIN") and leave the
# five new placeholder columns (C:G) blank for now - they'll be populated by
# the upcoming Monte Carlo simulation.
$ws.Range("C2:G21").ClearContents()

# Team / ExpPoints data, final state (row order + refreshed values).
$data = @(
    @{ Row = 2;  Team = "Real Madrid";         Exp = 88.95159428044995 },
    @{ Row = 3;  Team = "Barcelona";            Exp = 88.55144180756697 },
    @{ Row = 4;  Team = "Atlético de Madrid";   Exp = 70.38232796630246 },
    @{ Row = 5;  Team = "Villarreal";           Exp = 64.26874369684572 },
    @{ Row = 6;  Team = "Real Betis";           Exp = 60.0366739851554 },
    @{ Row = 7;  Team = "Athletic Club";        Exp = 59.78065867930323 },
    @{ Row = 8;  Team = "Rayo Vallecano";       Exp = 56.12249532014892 },
    @{ Row = 9;  Team = "Osasuna";              Exp = 49.28248075841894 },
    @{ Row = 10; Team = "Valencia";             Exp = 48.54234922738399 },
    @{ Row = 11; Team = "Sevilla";              Exp = 47.88776615197553 },
    @{ Row = 12; Team = "Real Sociedad";        Exp = 47.6004195287753 },
    @{ Row = 13; Team = "Getafe";               Exp = 47.49140136194507 },
    @{ Row = 14; Team = "Espanyol";             Exp = 47.15676655914623 },
    @{ Row = 15; Team = "Celta de Vigo";        Exp = 46.92979146006359 },
    @{ Row = 16; Team = "Mallorca";             Exp = 40.01750000886474 },
    @{ Row = 17; Team = "Alavés";               Exp = 38.6623134426073 },
    @{ Row = 18; Team = "Elche";                Exp = 36.12771650165123 },
    @{ Row = 19; Team = "Girona";               Exp = 34.26223879124277 },
    @{ Row = 20; Team = "Levante";              Exp = 33.26169983840459 },
    @{ Row = 21; Team = "Real Oviedo";          Exp = 28.59535783338381 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.Team
    $ws.Cells.Item($item.Row, 8).Value = $item.Exp
}

$ws.Range("A1").Select()
